$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The template had several "dummy"/leftover sample rows (the automation used
# to append one row per run). When the automation finishes, those extra
# sample rows are removed, leaving only the real registered rows.
# Remove rows 9-11 first (bottom-up) so row numbers for the next delete stay valid.
$ws.Rows("9:11").Delete()
$ws.Rows("4:6").Delete()

# Leave the window in the state it's in right after the automation finishes:
# selection parked on C8.
$ws.Range("C8").Select()
